# "Login related testcases without validation"
# - fix casing of the username header ("userName" -> "username")
# - bump the sheet's internal sheetId (1 -> 2) by cloning the sheet and
#   dropping the original, keeping the same name/position/data
# - move the saved cursor selection to C5
# - tighten the column widths (approximating Excel's AutoFit pass)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text casing.
$ws.Range("A1").Value = "username"

# Re-create the sheet so it gets a fresh (higher) sheetId, same as the
# original author's save did. Copying preserves data/styles/selection;
# we then delete the old sheet and rename the clone back. (References
# captured before the Delete() go stale once indices shift, so the
# clone is re-fetched by name afterwards rather than reusing a handle.)
$origName = $ws.Name
$cloneName = $origName + " (2)"
[void]$ws.Copy($null, $ws)
[void]$wb.Worksheets($origName).Delete()
$wb.Worksheets($cloneName).Name = $origName
[void]$wb.Worksheets($origName).Activate()

$ws = $wb.ActiveSheet

# Update the saved selection.
[void]$ws.Range("C5").Select()

# Approximate the AutoFit column-width pass.
$ws.Columns("A").ColumnWidth = 17.666666666666668
$ws.Columns("B").ColumnWidth = 10.998697916666666
